# Update "想去人数" (number of people wanting to go) counts in the
# "展览" sheet and the mirrored "全部类型" aggregate sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 11733
$ws1.Range("F8").Value = 4398
$ws1.Range("F13").Value = 2551
$ws1.Range("F20").Value = 519
$ws1.Range("F22").Value = 11290

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 11733
$ws4.Range("F8").Value = 4398
$ws4.Range("F13").Value = 2551
$ws4.Range("F21").Value = 519
$ws4.Range("F23").Value = 11290
